$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows 2-52 down to 3-53
$ws.Rows.Item(2).Insert()

# The inserted row inherits header-like formatting; clear it so it starts from default (style 0)
$ws.Range("A2:E2").ClearFormats()

# Populate new row 2 with the new data point (2007/2008 pair)
$ws.Range("A2").Value() = 39400
$ws.Range("B2").Value() = 2007
$ws.Range("C2").Value() = 11.13090654781819
$ws.Range("D2").Value() = 2008
$ws.Range("E2").Value() = 10.67037004222142

# Re-apply the date-style formatting (style index 2 in the original workbook) to A2,
# matching the rest of column A (bold font, thin border, centered/top alignment, date number format)
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").Borders.LineStyle = 1

# Update the remaining (shifted) rows 3-53 with the recomputed values from the diff
# Row 3
$ws.Range("A3").Value() = 39583
$ws.Range("B3").Value() = 2008
$ws.Range("C3").Value() = 7.288845540350142
$ws.Range("D3").Value() = 2009
$ws.Range("E3").Value() = 10.06916370210014

# Row 4
$ws.Range("A4").Value() = 39765
$ws.Range("B4").Value() = 2008
$ws.Range("C4").Value() = 4.672550446571067
$ws.Range("D4").Value() = 2009
$ws.Range("E4").Value() = -0.7156496512470745

# Row 5
$ws.Range("A5").Value() = 39948
$ws.Range("B5").Value() = 2009
$ws.Range("C5").Value() = -20.40984652067478
$ws.Range("D5").Value() = 2010
$ws.Range("E5").Value() = -33.45158075171641

# Row 6
$ws.Range("A6").Value() = 40130
$ws.Range("B6").Value() = 2009
$ws.Range("C6").Value() = -14.45332333832743
$ws.Range("D6").Value() = 2010
$ws.Range("E6").Value() = 7.857938327064184

# Row 7
$ws.Range("A7").Value() = 40310
$ws.Range("B7").Value() = 2010
$ws.Range("C7").Value() = 5.331710924091837
$ws.Range("D7").Value() = 2011
$ws.Range("E7").Value() = 10.84949786623359

# Row 8
$ws.Range("A8").Value() = 40494
$ws.Range("B8").Value() = 2010
$ws.Range("C8").Value() = 8.600536527919633
$ws.Range("D8").Value() = 2011
$ws.Range("E8").Value() = 12.64892828543749

# Row 9
$ws.Range("A9").Value() = 40676
$ws.Range("B9").Value() = 2011
$ws.Range("C9").Value() = 11.04982736891558
$ws.Range("D9").Value() = 2012
$ws.Range("E9").Value() = 9.589921161142879

# Row 10
$ws.Range("A10").Value() = 40862
$ws.Range("B10").Value() = 2011
$ws.Range("C10").Value() = 10.25770250047622
$ws.Range("D10").Value() = 2012
$ws.Range("E10").Value() = 10.40099841437159

# Row 11
$ws.Range("A11").Value() = 41044
$ws.Range("B11").Value() = 2012
$ws.Range("C11").Value() = 4.748210439985256
$ws.Range("D11").Value() = 2013
$ws.Range("E11").Value() = 7.086193663490992

# Row 12
$ws.Range("A12").Value() = 41228
$ws.Range("B12").Value() = 2012
$ws.Range("C12").Value() = 4.639893381363169
$ws.Range("D12").Value() = 2013
$ws.Range("E12").Value() = 6.662398279632087

# Row 13
$ws.Range("A13").Value() = 41409
$ws.Range("B13").Value() = 2013
$ws.Range("C13").Value() = -2.313034291448757
$ws.Range("D13").Value() = 2014
$ws.Range("E13").Value() = -7.020874871669158

# Row 14
$ws.Range("A14").Value() = 41592
$ws.Range("B14").Value() = 2013
$ws.Range("C14").Value() = 0.3058963467304165
$ws.Range("D14").Value() = 2014
$ws.Range("E14").Value() = 1.195213983078647

# Row 15
$ws.Range("A15").Value() = 41774
$ws.Range("B15").Value() = 2014
$ws.Range("C15").Value() = 3.711391384148
$ws.Range("D15").Value() = 2015
$ws.Range("E15").Value() = 0.6610955960690834

# Row 16
$ws.Range("A16").Value() = 41957
$ws.Range("B16").Value() = 2014
$ws.Range("C16").Value() = 4.068173739091874
$ws.Range("D16").Value() = 2015
$ws.Range("E16").Value() = 7.055025120039615

# Row 17
$ws.Range("A17").Value() = 42137
$ws.Range("B17").Value() = 2015
$ws.Range("C17").Value() = 3.860244074450181
$ws.Range("D17").Value() = 2016
$ws.Range("E17").Value() = 3.254220449867029

# Row 18
$ws.Range("A18").Value() = 42321
$ws.Range("B18").Value() = 2015
$ws.Range("C18").Value() = 4.984288257750213
$ws.Range("D18").Value() = 2016
$ws.Range("E18").Value() = 1.985659800779893

# Row 19
$ws.Range("A19").Value() = 42503
$ws.Range("B19").Value() = 2016
$ws.Range("C19").Value() = 2.370939381494686
$ws.Range("D19").Value() = 2017
$ws.Range("E19").Value() = 3.967543131851214

# Row 20
$ws.Range("A20").Value() = 42689
$ws.Range("B20").Value() = 2016
$ws.Range("C20").Value() = 1.878184267712912
$ws.Range("D20").Value() = 2017
$ws.Range("E20").Value() = -0.3562142672005275

# Row 21
$ws.Range("A21").Value() = 42867
$ws.Range("B21").Value() = 2017
$ws.Range("C21").Value() = 4.421855465610292
$ws.Range("D21").Value() = 2018
$ws.Range("E21").Value() = 5.260364862099642

# Row 22
$ws.Range("A22").Value() = 43053
$ws.Range("B22").Value() = 2017
$ws.Range("C22").Value() = 4.695933104194339
$ws.Range("D22").Value() = 2018
$ws.Range("E22").Value() = 6.493919935864634

# Row 23
$ws.Range("A23").Value() = 43145
$ws.Range("B23").Value() = 2018
$ws.Range("C23").Value() = 9.429264335267163
$ws.Range("D23").Value() = 2019
$ws.Range("E23").Value() = 11.40563592910597

# Row 24
$ws.Range("A24").Value() = 43235
$ws.Range("B24").Value() = 2018
$ws.Range("C24").Value() = 3.320585727896552
$ws.Range("D24").Value() = 2019
$ws.Range("E24").Value() = -4.098213472638578

# Row 25
$ws.Range("A25").Value() = 43326
$ws.Range("B25").Value() = 2018
$ws.Range("C25").Value() = 5.276665321936447
$ws.Range("D25").Value() = 2019
$ws.Range("E25").Value() = 2.51398037094428

# Row 26
$ws.Range("A26").Value() = 43418
$ws.Range("B26").Value() = 2018
$ws.Range("C26").Value() = 4.892602738886098
$ws.Range("D26").Value() = 2019
$ws.Range("E26").Value() = -2.576675125869599

# Row 27
$ws.Range("A27").Value() = 43510
$ws.Range("B27").Value() = 2019
$ws.Range("C27").Value() = 1.320842979722947
$ws.Range("D27").Value() = 2020
$ws.Range("E27").Value() = 2.636028935395296

# Row 28
$ws.Range("A28").Value() = 43600
$ws.Range("B28").Value() = 2019
$ws.Range("C28").Value() = 1.782333336406405
$ws.Range("D28").Value() = 2020
$ws.Range("E28").Value() = 4.060401000000002

# Row 29
$ws.Range("A29").Value() = 43691
$ws.Range("B29").Value() = 2019
$ws.Range("C29").Value() = 0.3252781783188663
$ws.Range("D29").Value() = 2020
$ws.Range("E29").Value() = -4.518236404743526

# Row 30
$ws.Range("A30").Value() = 43783
$ws.Range("B30").Value() = 2019
$ws.Range("C30").Value() = 0.8049382522247184
$ws.Range("D30").Value() = 2020
$ws.Range("E30").Value() = 3.1919852842623

# Row 31
$ws.Range("A31").Value() = 43875
$ws.Range("B31").Value() = 2020
$ws.Range("C31").Value() = -0.1289008616491616
$ws.Range("D31").Value() = 2021
$ws.Range("E31").Value() = -0.9308772335758553

# Row 32
$ws.Range("A32").Value() = 43966
$ws.Range("B32").Value() = 2020
$ws.Range("C32").Value() = -4.477718018907028
$ws.Range("D32").Value() = 2021
$ws.Range("E32").Value() = -11.83522404790002

# Row 33
$ws.Range("A33").Value() = 44068
$ws.Range("B33").Value() = 2020
$ws.Range("C33").Value() = -9.810777096850787
$ws.Range("D33").Value() = 2021
$ws.Range("E33").Value() = 29.96709940045137

# Row 34
$ws.Range("A34").Value() = 44159
$ws.Range("B34").Value() = 2020
$ws.Range("C34").Value() = -8.784173899737169
$ws.Range("D34").Value() = 2021
$ws.Range("E34").Value() = 6.942816049735523

# Row 35
$ws.Range("A35").Value() = 44251
$ws.Range("B35").Value() = 2021
$ws.Range("C35").Value() = 7.026336004273714
$ws.Range("D35").Value() = 2022
$ws.Range("E35").Value() = 7.425901647531985

# Row 36
$ws.Range("A36").Value() = 44341
$ws.Range("B36").Value() = 2021
$ws.Range("C36").Value() = 6.317691071509768
$ws.Range("D36").Value() = 2022
$ws.Range("E36").Value() = 4.613033063261129

# Row 37
$ws.Range("A37").Value() = 44432
$ws.Range("B37").Value() = 2021
$ws.Range("C37").Value() = 5.797134106720514
$ws.Range("D37").Value() = 2022
$ws.Range("E37").Value() = 5.455672087096408

# Row 38
$ws.Range("A38").Value() = 44525
$ws.Range("B38").Value() = 2021
$ws.Range("C38").Value() = 5.110501195359984
$ws.Range("D38").Value() = 2022
$ws.Range("E38").Value() = 0.8094958705429534

# Row 39
$ws.Range("A39").Value() = 44617
$ws.Range("B39").Value() = 2022
$ws.Range("C39").Value() = 2.670821531651923
$ws.Range("D39").Value() = 2023
$ws.Range("E39").Value() = -7.844778655777695

# Row 40
$ws.Range("A40").Value() = 44706
$ws.Range("B40").Value() = 2022
$ws.Range("C40").Value() = 3.79744344971964
$ws.Range("D40").Value() = 2023
$ws.Range("E40").Value() = 3.703837953294542

# Row 41
$ws.Range("A41").Value() = 44798
$ws.Range("B41").Value() = 2022
$ws.Range("C41").Value() = 4.232564748995715
$ws.Range("D41").Value() = 2023
$ws.Range("E41").Value() = 0.3861805562020093

# Row 42
$ws.Range("A42").Value() = 44890
$ws.Range("B42").Value() = 2022
$ws.Range("C42").Value() = 5.120680133083599
$ws.Range("D42").Value() = 2023
$ws.Range("E42").Value() = 0.5542886326586061

# Row 43
$ws.Range("A43").Value() = 44981
$ws.Range("B43").Value() = 2023
$ws.Range("C43").Value() = -0.2545313393182314
$ws.Range("D43").Value() = 2024
$ws.Range("E43").Value() = -2.527634545037938

# Row 44
$ws.Range("A44").Value() = 45071
$ws.Range("B44").Value() = 2023
$ws.Range("C44").Value() = 0.4998689793225486
$ws.Range("D44").Value() = 2024
$ws.Range("E44").Value() = -0.08221002454066317

# Row 45
$ws.Range("A45").Value() = 45163
$ws.Range("B45").Value() = 2023
$ws.Range("C45").Value() = 0.08070151925247959
$ws.Range("D45").Value() = 2024
$ws.Range("E45").Value() = 0.7505831475431046

# Row 46
$ws.Range("A46").Value() = 45254
$ws.Range("B46").Value() = 2023
$ws.Range("C46").Value() = -0.5532735011319234
$ws.Range("D46").Value() = 2024
$ws.Range("E46").Value() = -3.561435976944571

# Row 47
$ws.Range("A47").Value() = 45345
$ws.Range("B47").Value() = 2024
$ws.Range("C47").Value() = -2.978154922642562
$ws.Range("D47").Value() = 2025
$ws.Range("E47").Value() = -1.015253537920036

# Row 48
$ws.Range("A48").Value() = 45436
$ws.Range("B48").Value() = 2024
$ws.Range("C48").Value() = 0.1967053802870877
$ws.Range("D48").Value() = 2025
$ws.Range("E48").Value() = 6.289039804796182

# Row 49
$ws.Range("A49").Value() = 45534
$ws.Range("B49").Value() = 2024
$ws.Range("C49").Value() = -0.9685570952743805
$ws.Range("D49").Value() = 2025
$ws.Range("E49").Value() = -0.5620920786801986

# Row 50
$ws.Range("A50").Value() = 45618
$ws.Range("B50").Value() = 2024
$ws.Range("C50").Value() = -1.069674659641462
$ws.Range("D50").Value() = 2025
$ws.Range("E50").Value() = 0.01743232028155184

# Row 51
$ws.Range("A51").Value() = 45713
$ws.Range("B51").Value() = 2025
$ws.Range("C51").Value() = -4.169154013177412
$ws.Range("D51").Value() = 2026
$ws.Range("E51").Value() = -2.751196593554839

# Row 52
$ws.Range("A52").Value() = 45800
$ws.Range("B52").Value() = 2025
$ws.Range("C52").Value() = -2.056549539789942
$ws.Range("D52").Value() = 2026
$ws.Range("E52").Value() = -0.3884660724497446

# Row 53
$ws.Range("A53").Value() = 45891
$ws.Range("B53").Value() = 2025
$ws.Range("C53").Value() = -2.436529450546909
$ws.Range("D53").Value() = 2026
$ws.Range("E53").Value() = 0.4756432387424292
